$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Content changes: the depreciated asset changed from a gaming PC to a
#    4K monitor, and the asset code / DNI values were updated.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Monitor"
$ws.Range("D3").Value = "vybDy162"

# E3 must stay a text value ("001") rather than being converted to the
# number 1, so force a text format before assigning it, then reset the
# cell style back to Normal so no stray number-format style lingers.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "001"
$ws.Range("E3").Style = "Normal"

$ws.Range("C10").Value = "Monitor 4K"
$ws.Range("C8").Value = 5000
$ws.Range("C9").Value = 3

# ---------------------------------------------------------------------------
# 2. Formatting changes: add a thin gray border around the header, label and
#    value cells, make the font color explicit (automatic/black) for the
#    header & labels, and switch the header fill from the dotted yellow
#    pattern to a solid yellow fill (matching the label cells).
# ---------------------------------------------------------------------------
$grayColor = 8421504   # RGB(128,128,128)
$yellowColor = 65535   # RGB(255,255,0)

# Header row B2:E2
$header = $ws.Range("B2")
$header.Borders.Color = $grayColor
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.Interior.Pattern = 1
$header.Interior.Color = $yellowColor
$header.Font.ThemeColor = 1
$header.Copy()
foreach ($addr in @("C2", "D2", "E2")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# Label cells B6, B8, B9, B10
$label = $ws.Range("B6")
$label.Borders.Color = $grayColor
$label.Borders.LineStyle = 1
$label.Borders.Weight = 2
$label.Font.ThemeColor = 1
$label.Copy()
foreach ($addr in @("B8", "B9", "B10")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# Value cells C8, C9, C10
$value = $ws.Range("C8")
$value.Borders.Color = $grayColor
$value.Borders.LineStyle = 1
$value.Borders.Weight = 2
$value.Copy()
foreach ($addr in @("C9", "C10")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Row 6 no longer needs an explicit custom row height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).AutoFit()
